# Adds Viktor's own sentiment evaluation in column D (header + per-sentence
# neg/poz/neut labels), mirroring the commit: "Dodao svoju evaluaciju".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Viktor"
$ws.Range("D2").Value = "Neg"
$ws.Range("D3").Value = "neg"
$ws.Range("D4").Value = "neut"
$ws.Range("D5").Value = "poz"
$ws.Range("D6").Value = "neg"
$ws.Range("D7").Value = "neg"
$ws.Range("D8").Value = "neg"
$ws.Range("D9").Value = "neg"
$ws.Range("D10").Value = "neut"
$ws.Range("D11").Value = "neut"
$ws.Range("D12").Value = "neut"
$ws.Range("D13").Value = "neut"
$ws.Range("D14").Value = "neut"
$ws.Range("D15").Value = "poz"
$ws.Range("D16").Value = "poz"
$ws.Range("D17").Value = "poz"
$ws.Range("D18").Value = "poz"
$ws.Range("D19").Value = "poz"
$ws.Range("D20").Value = "poz"
$ws.Range("D21").Value = "neg"
$ws.Range("D22").Value = "poz"
$ws.Range("D23").Value = "poz"
$ws.Range("D24").Value = "poz"
$ws.Range("D25").Value = "neg"
$ws.Range("D27").Value = "neut"
$ws.Range("D28").Value = "neut"
$ws.Range("D29").Value = "neut"
$ws.Range("D30").Value = "poz"
$ws.Range("D31").Value = "neut"
$ws.Range("D32").Value = "neut"
$ws.Range("D33").Value = "poz"
$ws.Range("D34").Value = "neut"
$ws.Range("D35").Value = "poz"
$ws.Range("D36").Value = "poz"
$ws.Range("D37").Value = "poz"
$ws.Range("D38").Value = "poz"
$ws.Range("D39").Value = "poz"
$ws.Range("D40").Value = "neut"
$ws.Range("D41").Value = "neg"
$ws.Range("D42").Value = "neut"
$ws.Range("D43").Value = "neut"
$ws.Range("D44").Value = "neut"
$ws.Range("D45").Value = "neut"
$ws.Range("D46").Value = "neut"
$ws.Range("D47").Value = "neut"
$ws.Range("D48").Value = "poz"
$ws.Range("D49").Value = "neut"
$ws.Range("D51").Value = "neg"
$ws.Range("D52").Value = "poz"
$ws.Range("D53").Value = "poz"
$ws.Range("D54").Value = "poz"
$ws.Range("D55").Value = "poz"
$ws.Range("D56").Value = "poz"
$ws.Range("D57").Value = "poz"
$ws.Range("D58").Value = "poz"
$ws.Range("D59").Value = "poz"
$ws.Range("D60").Value = "poz"
$ws.Range("D61").Value = "neut"
$ws.Range("D62").Value = "neut"
$ws.Range("D63").Value = "poz"
$ws.Range("D64").Value = "poz"
$ws.Range("D65").Value = "poz"
$ws.Range("D66").Value = "poz"
$ws.Range("D67").Value = "neg"
$ws.Range("D68").Value = "poz"
$ws.Range("D69").Value = "neg"
$ws.Range("D70").Value = "poz"
$ws.Range("D71").Value = "poz"
$ws.Range("D72").Value = "poz"
$ws.Range("D73").Value = "neut"
$ws.Range("D74").Value = "poz"
$ws.Range("D75").Value = "poz"
$ws.Range("D76").Value = "neut"
$ws.Range("D77").Value = "poz"
$ws.Range("D78").Value = "neg"
$ws.Range("D80").Value = "neut"
$ws.Range("D81").Value = "poz"
$ws.Range("D82").Value = "neg"
$ws.Range("D83").Value = "neg"
$ws.Range("D84").Value = "poz"
$ws.Range("D85").Value = "poz"
$ws.Range("D86").Value = "neg"
$ws.Range("D87").Value = "poz"
$ws.Range("D88").Value = "poz"
$ws.Range("D89").Value = "poz"
$ws.Range("D90").Value = "poz"
$ws.Range("D91").Value = "neg"
$ws.Range("D92").Value = "neg"
$ws.Range("D93").Value = "neg"
$ws.Range("D94").Value = "neut"
$ws.Range("D95").Value = "neut"
$ws.Range("D96").Value = "poz"
$ws.Range("D97").Value = "neut"
$ws.Range("D98").Value = "neut"
$ws.Range("D99").Value = "poz"
$ws.Range("D100").Value = "poz"
$ws.Range("D101").Value = "neut"
$ws.Range("D102").Value = "neut"
$ws.Range("D103").Value = "neg"
$ws.Range("D104").Value = "neg"
$ws.Range("D105").Value = "neg"
$ws.Range("D106").Value = "neg"
$ws.Range("D107").Value = "neut"
$ws.Range("D109").Value = "poz"
$ws.Range("D110").Value = "poz"
$ws.Range("D111").Value = "poz"
$ws.Range("D112").Value = "poz"
$ws.Range("D113").Value = "poz"
$ws.Range("D114").Value = "poz"
$ws.Range("D115").Value = "poz"
$ws.Range("D116").Value = "poz"
$ws.Range("D117").Value = "neg"
$ws.Range("D118").Value = "neg"
$ws.Range("D119").Value = "poz"
$ws.Range("D120").Value = "poz"
$ws.Range("D121").Value = "neg"
$ws.Range("D122").Value = "neut"
$ws.Range("D123").Value = "neg"
$ws.Range("D124").Value = "neut"
$ws.Range("D125").Value = "neut"
$ws.Range("D126").Value = "neut"
$ws.Range("D127").Value = "neut"
$ws.Range("D128").Value = "poz"
$ws.Range("D129").Value = "poz"
$ws.Range("D130").Value = "neg"
$ws.Range("D131").Value = "poz"
$ws.Range("D132").Value = "poz"
$ws.Range("D133").Value = "neut"
$ws.Range("D135").Value = "poz"
$ws.Range("D136").Value = "poz"
$ws.Range("D137").Value = "poz"
$ws.Range("D138").Value = "poz"
$ws.Range("D139").Value = "poz"
$ws.Range("D140").Value = "neuz"
$ws.Range("D141").Value = "poz"
$ws.Range("D142").Value = "poz"
$ws.Range("D143").Value = "poz"
$ws.Range("D144").Value = "neut"
$ws.Range("D145").Value = "neut"
$ws.Range("D146").Value = "poz"
$ws.Range("D147").Value = "neut"
$ws.Range("D148").Value = "poz"
$ws.Range("D149").Value = "poz"
$ws.Range("D150").Value = "poz"
$ws.Range("D151").Value = "poz"
$ws.Range("D152").Value = "poz"

# Best-effort view-state bookkeeping (selection / zoom) to mirror the
# author's last on-screen state when they saved the file.
$ws.Range("D152").Select()
$excel.ActiveWindow.Zoom = 112
$excel.ActiveWindow.ScrollRow = 125
$excel.ActiveWindow.ScrollColumn = 1
